$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before row 610, shifting existing rows (and the
# rest of the weekly Acelga price history) down by two rows. This mirrors the
# author inserting a brand-new week's worth of records (dated 2022-01-17 /
# serial 44578) at the top of the "Mercado Mayorista Lo Valledor de Santiago"
# block.
$ws.Range("A610:A611").EntireRow.Insert()

# --- New row 610: Acelga, Primera, week of 2022-01-17 ---
$ws.Range("A610").Value = 6
$ws.Range("B610").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C610").Value = "Metropolitana"
$ws.Range("D610").Value2 = 44578
$ws.Range("E610").Value = 13
$ws.Range("F610").Value = 100112009
$ws.Range("G610").Value = "Acelga"
$ws.Range("H610").Value = "Sin especificar"
$ws.Range("I610").Value = "Primera"
$ws.Range("J610").Value = 130
$ws.Range("K610").Value = 16000
$ws.Range("L610").Value = 16000
$ws.Range("M610").Value = 16000
$ws.Range("N610").Value = "$/docena de atados"
$ws.Range("O610").Value = "Región Metropolitana"
$ws.Range("P610").Value = 5333
$ws.Range("Q610").Value = 3
$ws.Range("R610").Value = "Hortaliza"

# --- New row 611: Acelga, Segunda, week of 2022-01-17 ---
$ws.Range("A611").Value = 6
$ws.Range("B611").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C611").Value = "Metropolitana"
$ws.Range("D611").Value2 = 44578
$ws.Range("E611").Value = 13
$ws.Range("F611").Value = 100112009
$ws.Range("G611").Value = "Acelga"
$ws.Range("H611").Value = "Sin especificar"
$ws.Range("I611").Value = "Segunda"
$ws.Range("J611").Value = 60
$ws.Range("K611").Value = 14000
$ws.Range("L611").Value = 14000
$ws.Range("M611").Value = 14000
$ws.Range("N611").Value = "$/docena de atados"
$ws.Range("O611").Value = "Región Metropolitana"
$ws.Range("P611").Value = 4667
$ws.Range("Q611").Value = 3
$ws.Range("R611").Value = "Hortaliza"

Write-Host "Inserted new week (44578) rows at 610-611"
